$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 18: fill in previously-blank cells with values ----
$ws.Range("C18").Value = 2010
$ws.Range("D18").Value = 1090.199341
$ws.Range("E18").Value = 1990.4676509999999
$ws.Range("F18").Value = 1.2021059999999999
$ws.Range("G18").Value = 280.16485599999999
$ws.Range("H18").Value = 10.610913999999999
$ws.Range("I18").Value = 4.714696
$ws.Range("J18").Value = 8.8404570000000007
$ws.Range("K18").Value = 677.33367899999996
$ws.Range("L18").Value = 93.229797000000005
$ws.Range("M18").Value = 1393.5979
$ws.Range("N18").Value = 1207.1701660000001
$ws.Range("O18").Value = 6798.6591799999997
$ws.Range("P18").Value = 29450.638672000001
$ws.Range("Q18").Value = 2.8124359999999999
$ws.Range("R18").Value = 0.00083299999999999997
$ws.Range("S18").Value = 2010

# R18 gets a tighter 5-decimal number format (new custom number format 0.00000)
$ws.Range("R18").NumberFormat = "0.00000"
# S18 was an empty formatted placeholder cell before; now it reverts to a plain,
# unformatted cell (General format) once it actually holds data
$ws.Range("S18").ClearFormats()

# ---- Row 30: new row, duplicate of row 29 (copy row 29 down one row) ----
$ws.Range("D29:R29").Copy()
$ws.Range("D30").PasteSpecial(-4122)   # xlPasteFormats - bring number formats along

$ws.Range("D30").Value = 1186.9540608888888
$ws.Range("E30").Value = 1901.5157334444443
$ws.Range("F30").Value = 0.97970299999999988
$ws.Range("G30").Value = 280.33542888888883
$ws.Range("H30").Value = 9.775355222222224
$ws.Range("I30").Value = 5.3768950000000002
$ws.Range("J30").Value = 8.145128999999999
$ws.Range("K30").Value = 645.94088411111113
$ws.Range("L30").Value = 83.47062044444445
$ws.Range("M30").Value = 1457.6001788888889
$ws.Range("N30").Value = 1191.1331380000001
$ws.Range("O30").Value = 4695.8937716666669
$ws.Range("P30").Value = 27227.338324888889
$ws.Range("Q30").Value = 1.3527740000000001
$ws.Range("R30").Value = 0.00039611111111111119

# The pasted-down row had its fill explicitly cleared afterwards
$ws.Range("D30:R30").Interior.ColorIndex = -4142

# Update the active selection to match the new last row
$ws.Range("A30:XFD30").Select()

Write-Host "done"
